$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append new row 57 with the new test-mail entry
$ws.Range("A57").Value = "Zou je dit kunnen doorsturen?"
$ws.Range("B57").Value = "mailmind.test@zohomail.eu"
$ws.Range("C57").Value = "Testmail #16: Zou je dit kunnen doorsturen?"
$ws.Range("D57").Value = "Overig"
$ws.Range("E57").Value = "Geachte klant,`nHartelijk dank voor uw e-mail. Om u beter van dienst te kunnen zijn, zou u ons meer informatie kunnen geven over wat u precies wilt laten doorsturen? Op die manier kunnen we u gerichter helpen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$ws.Range("F57").Value = "2025-08-05 19:56:03"
$ws.Range("G57").Value = "Ja"
$ws.Range("H57").Value = "Nee"
$ws.Range("I57").Value = "Ja"
$ws.Range("J57").Value = "Nee"

# The multi-line text in column E auto-expands row 57's height; re-run
# AutoFit (with no explicit WrapText) so it settles back on the sheet's
# default row height, matching every other data row.
$ws.Rows.Item(57).AutoFit()

# Extend the conditional-formatting blocks (one per column) so they keep
# covering the data range through the newly added row.
$ws.Range("D2:D56").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D57"))
$ws.Range("G2:G56").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G57"))
$ws.Range("H2:H56").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H57"))
$ws.Range("I2:I56").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I57"))
$ws.Range("J2:J56").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J57"))

# Update the Dashboard summary count for the "Overig" category.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 11
